$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.07285266666666666
$ws.Range("H2").Value = 0.218558
$ws.Range("I2").Value = 0.05584899373277382
$ws.Range("J2").Value = 0.05584899373277381
$ws.Range("M2").Value = 0.032838
$ws.Range("N2").Value = 0.098514
$ws.Range("O2").Value = 0.007146324094219707
$ws.Range("P2").Value = 0.007146324094219707
$ws.Range("Q2").Value = 0.002392335868
$ws.Range("R2").Value = 0.021531022812
$ws.Range("S2").Value = 0.0003991150095504469
$ws.Range("T2").Value = 0.0003991150095504469
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.07285266666666666
$ws.Range("H3").Value = 0.218558
$ws.Range("I3").Value = 0.05584899373277382
$ws.Range("J3").Value = 0.05584899373277381
$ws.Range("O3").Value = 0.03951718316124263
$ws.Range("P3").Value = 0.03951718316124263
$ws.Range("Q3").Value = 0.01322895147666667
$ws.Range("R3").Value = 0.11906056329
$ws.Range("S3").Value = 0.002206994914709115
$ws.Range("T3").Value = 0.002206994914709114
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.07285266666666666
$ws.Range("H4").Value = 0.218558
$ws.Range("I4").Value = 0.05584899373277382
$ws.Range("J4").Value = 0.05584899373277381
$ws.Range("M4").Value = 3.814633
$ws.Range("N4").Value = 11.443899
$ws.Range("O4").Value = 0.8301542030119253
$ws.Range("P4").Value = 0.8301542030119253
$ws.Range("Q4").Value = 0.2779061864046667
$ws.Range("R4").Value = 2.501155677642
$ws.Range("S4").Value = 0.04636327688124886
$ws.Range("T4").Value = 0.04636327688124886
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.07285266666666666
$ws.Range("H5").Value = 0.218558
$ws.Range("I5").Value = 0.05584899373277382
$ws.Range("J5").Value = 0.05584899373277381
$ws.Range("M5").Value = 0.5660336666666667
$ws.Range("N5").Value = 1.698101
$ws.Range("O5").Value = 0.1231822897326124
$ws.Range("P5").Value = 0.1231822897326124
$ws.Range("Q5").Value = 0.04123706203977778
$ws.Range("R5").Value = 0.371133558358
$ws.Range("S5").Value = 0.006879606927265399
$ws.Range("T5").Value = 0.006879606927265397
$ws.Range("I6").Value = 0.438974399073536
$ws.Range("J6").Value = 0.438974399073536
$ws.Range("M6").Value = 0.032838
$ws.Range("N6").Value = 0.098514
$ws.Range("O6").Value = 0.007146324094219707
$ws.Range("P6").Value = 0.007146324094219707
$ws.Range("Q6").Value = 0.018803815966
$ws.Range("R6").Value = 0.169234343694
$ws.Range("S6").Value = 0.003137053324844827
$ws.Range("T6").Value = 0.003137053324844828
$ws.Range("I7").Value = 0.438974399073536
$ws.Range("J7").Value = 0.438974399073536
$ws.Range("O7").Value = 0.03951718316124263
$ws.Range("P7").Value = 0.03951718316124263
$ws.Range("S7").Value = 0.01734703173128534
$ws.Range("T7").Value = 0.01734703173128534
$ws.Range("I8").Value = 0.438974399073536
$ws.Range("J8").Value = 0.438974399073536
$ws.Range("M8").Value = 3.814633
$ws.Range("N8").Value = 11.443899
$ws.Range("O8").Value = 0.8301542030119253
$ws.Range("P8").Value = 0.8301542030119253
$ws.Range("Q8").Value = 2.184349135447667
$ws.Range("R8").Value = 19.659142219029
$ws.Range("S8").Value = 0.3644164424055302
$ws.Range("T8").Value = 0.3644164424055302
$ws.Range("I9").Value = 0.438974399073536
$ws.Range("J9").Value = 0.438974399073536
$ws.Range("M9").Value = 0.5660336666666667
$ws.Range("N9").Value = 1.698101
$ws.Range("O9").Value = 0.1231822897326124
$ws.Range("P9").Value = 0.1231822897326124
$ws.Range("Q9").Value = 0.3241242736634444
$ws.Range("R9").Value = 2.917118462971
$ws.Range("S9").Value = 0.05407387161187573
$ws.Range("T9").Value = 0.05407387161187573
$ws.Range("G10").Value = 0.594248
$ws.Range("H10").Value = 1.782744
$ws.Range("I10").Value = 0.4555516544035914
$ws.Range("J10").Value = 0.4555516544035914
$ws.Range("M10").Value = 0.032838
$ws.Range("N10").Value = 0.098514
$ws.Range("O10").Value = 0.007146324094219707
$ws.Range("P10").Value = 0.007146324094219707
$ws.Range("Q10").Value = 0.019513915824
$ws.Range("R10").Value = 0.175625242416
$ws.Range("S10").Value = 0.003255519764026034
$ws.Range("T10").Value = 0.003255519764026034
$ws.Range("G11").Value = 0.594248
$ws.Range("H11").Value = 1.782744
$ws.Range("I11").Value = 0.4555516544035914
$ws.Range("J11").Value = 0.4555516544035914
$ws.Range("O11").Value = 0.03951718316124263
$ws.Range("P11").Value = 0.03951718316124263
$ws.Range("Q11").Value = 0.10790652308
$ws.Range("R11").Value = 0.97115870772
$ws.Range("S11").Value = 0.01800211816647382
$ws.Range("T11").Value = 0.01800211816647382
$ws.Range("G12").Value = 0.594248
$ws.Range("H12").Value = 1.782744
$ws.Range("I12").Value = 0.4555516544035914
$ws.Range("J12").Value = 0.4555516544035914
$ws.Range("M12").Value = 3.814633
$ws.Range("N12").Value = 11.443899
$ws.Range("O12").Value = 0.8301542030119253
$ws.Range("P12").Value = 0.8301542030119253
$ws.Range("Q12").Value = 2.266838030984
$ws.Range("R12").Value = 20.401542278856
$ws.Range("S12").Value = 0.3781781205921775
$ws.Range("T12").Value = 0.3781781205921775
$ws.Range("G13").Value = 0.594248
$ws.Range("H13").Value = 1.782744
$ws.Range("I13").Value = 0.4555516544035914
$ws.Range("J13").Value = 0.4555516544035914
$ws.Range("M13").Value = 0.5660336666666667
$ws.Range("N13").Value = 1.698101
$ws.Range("O13").Value = 0.1231822897326124
$ws.Range("P13").Value = 0.1231822897326124
$ws.Range("Q13").Value = 0.3363643743493333
$ws.Range("R13").Value = 3.027279369144
$ws.Range("S13").Value = 0.05611589588091411
$ws.Range("T13").Value = 0.0561158958809141
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.06473366666666668
$ws.Range("H14").Value = 0.194201
$ws.Range("I14").Value = 0.04962495279009878
$ws.Range("J14").Value = 0.04962495279009878
$ws.Range("M14").Value = 0.032838
$ws.Range("N14").Value = 0.098514
$ws.Range("O14").Value = 0.007146324094219707
$ws.Range("P14").Value = 0.007146324094219707
$ws.Range("Q14").Value = 0.002125724146
$ws.Range("R14").Value = 0.019131517314
$ws.Range("S14").Value = 0.0003546359957983984
$ws.Range("T14").Value = 0.0003546359957983984
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.06473366666666668
$ws.Range("H15").Value = 0.194201
$ws.Range("I15").Value = 0.04962495279009878
$ws.Range("J15").Value = 0.04962495279009878
$ws.Range("O15").Value = 0.03951718316124263
$ws.Range("P15").Value = 0.03951718316124263
$ws.Range("Q15").Value = 0.01175466286166667
$ws.Range("R15").Value = 0.105791965755
$ws.Range("S15").Value = 0.001961038348774352
$ws.Range("T15").Value = 0.001961038348774352
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.06473366666666668
$ws.Range("H16").Value = 0.194201
$ws.Range("I16").Value = 0.04962495279009878
$ws.Range("J16").Value = 0.04962495279009878
$ws.Range("M16").Value = 3.814633
$ws.Range("N16").Value = 11.443899
$ws.Range("O16").Value = 0.8301542030119253
$ws.Range("P16").Value = 0.8301542030119253
$ws.Range("Q16").Value = 0.2469351810776667
$ws.Range("R16").Value = 2.222416629699
$ws.Range("S16").Value = 0.04119636313296888
$ws.Range("T16").Value = 0.04119636313296887
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.06473366666666668
$ws.Range("H17").Value = 0.194201
$ws.Range("I17").Value = 0.04962495279009878
$ws.Range("J17").Value = 0.04962495279009878
$ws.Range("M17").Value = 0.5660336666666667
$ws.Range("N17").Value = 1.698101
$ws.Range("O17").Value = 0.1231822897326124
$ws.Range("P17").Value = 0.1231822897326124
$ws.Range("Q17").Value = 0.03664143470011112
$ws.Range("R17").Value = 0.329772912301
$ws.Range("S17").Value = 0.006112915312557161
$ws.Range("T17").Value = 0.006112915312557159
